# char_frac.xlsx v0.2.3 - add UTF32 w/o BOM and ValidUtf8
# Fill in the KOI8 (J/K) columns for rows 14-17 and restyle the
# Windows-1251 (H/I) columns on those same rows to match the rest
# of the table (rows 3-13), then update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Set-CellValue {
    param($range, $value)
    $ws.Range($range).Value = $value
}

function Copy-CellFormat {
    param($sourceRange, $destRange)
    $ws.Range($sourceRange).Copy() | Out-Null
    $ws.Range($destRange).PasteSpecial($xlPasteFormats) | Out-Null
}

# --- Cell text values (existing shared strings, reused not duplicated) ---
Set-CellValue "H14" "0xEC"
Set-CellValue "I14" "0xCC"
Set-CellValue "J14" "0xCD"
Set-CellValue "K14" "0xED"

Set-CellValue "H15" "0xE4"
Set-CellValue "I15" "0xC4"
Set-CellValue "J15" "0xC4"
Set-CellValue "K15" "0xE4"

Set-CellValue "H16" "0xEF"
Set-CellValue "I16" "0xCF"
Set-CellValue "J16" "0xD0"
Set-CellValue "K16" "0xF0"

Set-CellValue "H17" "0xF3"
Set-CellValue "I17" "0xD3"
Set-CellValue "J17" "0xD5"
Set-CellValue "K17" "0xF5"

# --- Formatting: copy the styles already used elsewhere in the table ---
# Row 14
Copy-CellFormat "H3" "H14"
Copy-CellFormat "H3" "I14"
Copy-CellFormat "J3" "J14"
Copy-CellFormat "K3" "K14"

# Row 15
Copy-CellFormat "H3" "H15"
Copy-CellFormat "H3" "I15"
Copy-CellFormat "J4" "J15"
Copy-CellFormat "K4" "K15"

# Row 16
Copy-CellFormat "H3" "H16"
Copy-CellFormat "H3" "I16"
Copy-CellFormat "J4" "J16"
Copy-CellFormat "K4" "K16"

# Row 17
Copy-CellFormat "L3" "H17"
Copy-CellFormat "L3" "I17"
Copy-CellFormat "J2" "J17"
Copy-CellFormat "J2" "K17"

$excel.CutCopyMode = $false

# --- Final selection state, as left by the author ---
$ws.Range("H17:I17").Select() | Out-Null
